# Updated RAD Test Cases and data to handle FEIN/SSN and Failures when
# Payment Applications are not deployed in QA2.
#
# Updates the "Date" column (B2:B4) on Sheet1 with new timestamps
# reflecting a fresh test execution.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "Fri Sep 29 11:39:28 EDT 2023"
$ws.Range("B3").Value = "Fri Sep 29 11:39:41 EDT 2023"
$ws.Range("B4").Value = "Fri Sep 29 11:39:55 EDT 2023"
